# Auto-generated edit script applying numeric updates described by the commit diff.
# Updates cell values across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets (Table_* ranges).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1537.1395
$ws.Range("I132").Value = 1462.6842
$ws.Range("J132").Value = 2103
$ws.Range("K132").Value = 4388.0526
$ws.Range("L132").Value = 6309
$ws.Range("M132").Value = -1858.0526
$ws.Range("N132").Value = -11369
$ws.Range("H138").Value = 3061.4055
$ws.Range("I138").Value = 1996.625
$ws.Range("J138").Value = 3355.138
$ws.Range("K138").Value = 5989.875
$ws.Range("L138").Value = 10065.414
$ws.Range("M138").Value = -849.875
$ws.Range("N138").Value = -20345.414
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 7692.154
$ws.Range("I110").Value = 6599.8
$ws.Range("K110").Value = 6599.8
$ws.Range("M110").Value = -4554.8
$ws.Range("H122").Value = 1468
$ws.Range("I122").Value = 1468
$ws.Range("K122").Value = 4404
$ws.Range("M122").Value = -1954
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2691.05
$ws.Range("I20").Value = 2337.6
$ws.Range("K20").Value = 2337.6
$ws.Range("M20").Value = -2090.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10117.654
$ws.Range("I58").Value = 5695.091
$ws.Range("K58").Value = 5695.091
$ws.Range("M58").Value = -5492.091
$ws.Range("H82").Value = 53500
$ws.Range("J82").Value = 53500
$ws.Range("L82").Value = 53500
$ws.Range("N82").Value = -54222
$ws.Range("H85").Value = 53500
$ws.Range("J85").Value = 53500
$ws.Range("L85").Value = 53500
$ws.Range("N85").Value = -55996
$ws.Range("H87").Value = 64999.5
$ws.Range("J87").Value = 64999.5
$ws.Range("L87").Value = 64999.5
$ws.Range("N87").Value = -67371.5
$ws.Range("H90").Value = 64999.5
$ws.Range("J90").Value = 64999.5
$ws.Range("L90").Value = 194998.5
$ws.Range("N90").Value = -206854.5
$ws.Range("H99").Value = 5065.8335
$ws.Range("I99").Value = 4679
$ws.Range("K99").Value = 4679
$ws.Range("M99").Value = -3181
$ws.Range("H126").Value = 5065.8335
$ws.Range("I126").Value = 4679
$ws.Range("K126").Value = 14037
$ws.Range("M126").Value = -11567
$ws.Range("H134").Value = 4803.7
$ws.Range("I134").Value = 2921.4614
$ws.Range("J134").Value = 8299.286
$ws.Range("K134").Value = 8764.3842
$ws.Range("L134").Value = 24897.858
$ws.Range("M134").Value = -6229.3842
$ws.Range("N134").Value = -29967.858
$ws.Range("H136").Value = 10117.654
$ws.Range("I136").Value = 5695.091
$ws.Range("K136").Value = 17085.273
$ws.Range("M136").Value = -14535.273
$ws.Range("H141").Value = 255194.44
$ws.Range("J141").Value = 255194.44
$ws.Range("L141").Value = 255194.44
$ws.Range("N141").Value = -265554.44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1429658.9
$ws.Range("J117").Value = 1667780.4
$ws.Range("L117").Value = 5003341.199999999
$ws.Range("N117").Value = -5010225.199999999
$ws.Range("H131").Value = 14708197
$ws.Range("J131").Value = 2893.2173
$ws.Range("L131").Value = 8679.651899999999
$ws.Range("N131").Value = -18759.6519
$ws.Range("H137").Value = 5245.3125
$ws.Range("I137").Value = 3383.111
$ws.Range("K137").Value = 10149.333
$ws.Range("M137").Value = -5049.332999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 12050.25
$ws.Range("J19").Value = 12050.25
$ws.Range("L19").Value = 12050.25
$ws.Range("N19").Value = -12626.25
$ws.Range("H20").Value = 7512676.5
$ws.Range("J20").Value = 17182.6
$ws.Range("L20").Value = 17182.6
$ws.Range("N20").Value = -17672.6
$ws.Range("H70").Value = 11550.571
$ws.Range("I70").Value = 999.5
$ws.Range("K70").Value = 999.5
$ws.Range("M70").Value = -729.5
$ws.Range("H73").Value = 11550.571
$ws.Range("I73").Value = 999.5
$ws.Range("K73").Value = 999.5
$ws.Range("M73").Value = -63.5
$ws.Range("H124").Value = 65390
$ws.Range("J124").Value = 65390
$ws.Range("L124").Value = 65390
$ws.Range("N124").Value = -75210
$ws.Range("H132").Value = 3964.8667
$ws.Range("I132").Value = 3955.7083
$ws.Range("K132").Value = 11867.1249
$ws.Range("M132").Value = -9337.124899999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5623
$ws.Range("I7").Value = 5449.1665
$ws.Range("K7").Value = 5449.1665
$ws.Range("M7").Value = -5337.1665
$ws.Range("H40").Value = 4085.6924
$ws.Range("I40").Value = 4045
$ws.Range("J40").Value = 4133.1665
$ws.Range("K40").Value = 4045
$ws.Range("L40").Value = 4133.1665
$ws.Range("M40").Value = -3909
$ws.Range("N40").Value = -4405.1665
$ws.Range("H55").Value = 71428650
$ws.Range("I55").Value = 83333400
$ws.Range("K55").Value = 83333400
$ws.Range("M55").Value = -83333227
$ws.Range("H68").Value = 7100.4287
$ws.Range("I68").Value = 3250
$ws.Range("J68").Value = 8640.6
$ws.Range("K68").Value = 3250
$ws.Range("L68").Value = 8640.6
$ws.Range("M68").Value = -2501
$ws.Range("N68").Value = -10138.6
$ws.Range("H71").Value = 7100.4287
$ws.Range("I71").Value = 3250
$ws.Range("J71").Value = 8640.6
$ws.Range("K71").Value = 16250
$ws.Range("L71").Value = 43203
$ws.Range("M71").Value = -12506
$ws.Range("N71").Value = -50691
$ws.Range("H76").Value = 43999.4
$ws.Range("J76").Value = 43999.4
$ws.Range("L76").Value = 43999.4
$ws.Range("N76").Value = -44675.4
$ws.Range("H79").Value = 43999.4
$ws.Range("J79").Value = 43999.4
$ws.Range("L79").Value = 43999.4
$ws.Range("N79").Value = -46339.4
$ws.Range("H126").Value = 5623
$ws.Range("I126").Value = 5449.1665
$ws.Range("K126").Value = 16347.4995
$ws.Range("M126").Value = -13877.4995
$ws.Range("H136").Value = 48786010
$ws.Range("I136").Value = 30308408
$ws.Range("J136").Value = 125006110
$ws.Range("K136").Value = 90925224
$ws.Range("L136").Value = 375018330
$ws.Range("M136").Value = -90922674
$ws.Range("N136").Value = -375023430
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5496
$ws.Range("I14").Value = 5494.6665
$ws.Range("J14").Value = 5500
$ws.Range("K14").Value = 5494.6665
$ws.Range("L14").Value = 5500
$ws.Range("M14").Value = -5326.6665
$ws.Range("N14").Value = -5836
$ws.Range("H15").Value = 6986.6665
$ws.Range("I15").Value = 6994
$ws.Range("J15").Value = 6983
$ws.Range("K15").Value = 6994
$ws.Range("L15").Value = 6983
$ws.Range("M15").Value = -6706
$ws.Range("N15").Value = -7559
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H126").Value = 7093.067
$ws.Range("I126").Value = 4789.1
$ws.Range("J126").Value = 11701
$ws.Range("K126").Value = 14367.3
$ws.Range("L126").Value = 35103
$ws.Range("M126").Value = -11897.3
$ws.Range("N126").Value = -40043
$ws.Range("H136").Value = 4469.5
$ws.Range("I136").Value = 4005.92
$ws.Range("J136").Value = 8332.666999999999
$ws.Range("K136").Value = 12017.76
$ws.Range("L136").Value = 24998.001
$ws.Range("M136").Value = -9467.76
$ws.Range("N136").Value = -30098.001

Write-Host "Applied 185 cell updates across 8 sheets."
